# Auto-committed on 2022/06/17 週五 17:40:09.87
#
# - DBS (sheet2): viewport scroll/selection moved from C8 to A5, and a new
#   lookup row (custNoEq / "CustNo =" / "CustNo ASC") was added under the
#   existing findL2079 row.
# - DBD (sheet1): the saved scroll position (topLeftCell) is cleared while
#   it stays the active/selected tab.

$wb = $excel.ActiveWorkbook

$dbd = $wb.Worksheets.Item("DBD")
$dbs = $wb.Worksheets.Item("DBS")

# --- DBS: add the new CustNo-equality lookup row ------------------------
$dbs.Activate()

$dbs.Range("B4").Value = "CustNo ="
$dbs.Range("C4").Value = "CustNo ASC"
$dbs.Range("A4").Value = "custNoEq"

# Match the updated on-sheet selection/scroll (was C8, now A5); the frozen
# header pane (ySplit=1, topLeftCell A2) is left untouched.
$dbs.Range("A5").Select() | Out-Null

# --- DBD: stays the active tab, but its remembered scroll offset resets -
$dbd.Activate() | Out-Null
